$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.735.61"
$ws.Range("E2").Value = "  +1.28%  "

# Row 3
$ws.Range("D3").Value = "2.303.04"
$ws.Range("E3").Value = "  +0.18%  "

# Row 4
$ws.Range("E4").Value = "  +0.23%  "

# Row 5
$ws.Range("D5").Value = "'316.58"
$ws.Range("E5").Value = "  -0.49%  "

# Row 6
$ws.Range("D6").Value = "'103.83"
$ws.Range("E6").Value = "  +0.65%  "

# Row 7
$ws.Range("E7").Value = "  -0.16%  "

# Row 8
$ws.Range("E8").Value = "  +0.30%  "

# Row 9
$ws.Range("D9").Value = "'0.606"
$ws.Range("E9").Value = "  -0.47%  "

# Row 10
$ws.Range("D10").Value = "'39.72"
$ws.Range("E10").Value = "  +0.30%  "

# Row 11
$ws.Range("D11").Value = "'0.0908"
$ws.Range("E11").Value = "  -0.16%  "

# Row 12
$ws.Range("D12").Value = "'8.55"
$ws.Range("E12").Value = "  +2.66%  "

# Row 13
$ws.Range("E13").Value = "  +1.32%  "

# Row 14
$ws.Range("E14").Value = "  +3.85%  "

# Row 15
$ws.Range("D15").Value = "'15.41"
$ws.Range("E15").Value = "  +0.60%  "

# Row 16
$ws.Range("D16").Value = "2.654.09"
$ws.Range("E16").Value = "  +0.21%  "

# Row 17
$ws.Range("D17").Value = "2.327.51"
$ws.Range("E17").Value = "  +1.55%  "

# Row 18
$ws.Range("D18").Value = "42.683.99"
$ws.Range("E18").Value = "  +1.02%  "

# Row 19
$ws.Range("D19").Value = "'7.58"
$ws.Range("E19").Value = "  +2.52%  "

# Row 20
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.0000106"
$ws.Range("E20").Value = "  +0.32%  "

# Row 21
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").Value = "'13.72"
$ws.Range("E21").Value = "  +23.04%  "

# Row 22
$ws.Range("D22").Value = "'74.03"
$ws.Range("E22").Value = "  +0.71%  "

# Row 23
$ws.Range("E23").Value = "  -2.93%  "

# Row 24
$ws.Range("D24").Value = "'267.12"
$ws.Range("E24").Value = "  -4.78%  "

# Row 25
$ws.Range("E25").Value = "  -1.07%  "

# Row 26
$ws.Range("E26").Value = "  +0.33%  "

# Row 27
$ws.Range("D27").Value = "'10.95"
$ws.Range("E27").Value = "  +0.86%  "

# Row 28
$ws.Range("D28").Value = "'2.35"
$ws.Range("E28").Value = "  -3.48%  "

# Row 29
$ws.Range("D29").Value = "'22.71"
$ws.Range("E29").Value = "  -1.05%  "

# Row 30
$ws.Range("D30").Value = "'6.61"
$ws.Range("E30").Value = "  +12.83%  "

# Row 31
$ws.Range("D31").Value = "'37.61"
$ws.Range("E31").Value = "  +4.48%  "

# Row 32
$ws.Range("D32").Value = "'165.67"
$ws.Range("E32").Value = "  +0.95%  "

# Row 33
$ws.Range("D33").Value = "'0.0885"

# Row 34
$ws.Range("E34").Value = "  -2.83%  "

# Row 35
$ws.Range("D35").Value = "'2.59"
$ws.Range("E35").Value = "  -0.83%  "

# Row 36
$ws.Range("E36").Value = "  -0.23%  "

# Row 37
$ws.Range("D37").Value = "'4.58"
$ws.Range("E37").Value = "  -0.17%  "

# Row 38
$ws.Range("D38").Value = "'0.0354"
$ws.Range("E38").Value = "  +1.60%  "

# Row 39
$ws.Range("D39").Value = "'3.74"
$ws.Range("E39").Value = "  -0.38%  "

# Row 40
$ws.Range("D40").Value = "'2.73"
$ws.Range("E40").Value = "  -1.53%  "

# Row 41
$ws.Range("D41").Value = "'1.60"
$ws.Range("E41").Value = "  +10.11%  "

# Row 42
$ws.Range("D42").Value = "'70.83"
$ws.Range("E42").Value = "  +2.24%  "

# Row 43
$ws.Range("D43").Value = "'96.31"
$ws.Range("E43").Value = "  -4.30%  "

# Row 44
$ws.Range("E44").Value = "  +0.95%  "

# Row 45
$ws.Range("E45").Value = "  +0.26%  "

# Row 46
$ws.Range("D46").Value = "'12.55"
$ws.Range("E46").Value = "  +4.59%  "

# Row 47
$ws.Range("D47").Value = "'117.18"
$ws.Range("E47").Value = "  +4.93%  "

# Row 48
$ws.Range("D48").Value = "'80.15"
$ws.Range("E48").Value = "  +4.22%  "

# Row 49
$ws.Range("D49").Value = "1.667.59"
$ws.Range("E49").Value = "  +4.25%  "

# Row 50
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "'5.31"
$ws.Range("E50").Value = "  +0.38%  "

# Row 51
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "'8.89"
$ws.Range("E51").Value = "  -0.65%  "
